$d = $word.ActiveDocument

# The document currently ends with an otherwise-empty list paragraph that
# only holds the "_GoBack" bookmark. We add three new Q&A pairs after the
# existing content; the first question goes into that bookmark paragraph,
# and the bookmark itself ends up on the very last (new) paragraph.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$listTemplate = $d.Paragraphs.Item(1).Range.ListFormat.ListTemplate

# --- Question 1: goes into the existing (bookmarked) trailing paragraph ---
$lastPara.Range.InsertBefore("What special folders are in the ASP.NET and what is their role?")

# --- Answer 1: new plain paragraph after Question 1 ---
$r = $d.Content
$r.Find.Execute("What special folders are in the ASP.NET and what is their role?")
$r.Collapse(0)
$r.InsertParagraphAfter()
$ans1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$ans1.Range.ListFormat.RemoveNumbers()
$ans1.Style = "Normal"
$ans1.Range.Text = "The special folders in ASP.NET is the Web.Config Folder, which is the configuration and settings file for a ASP.Net web application and contains data about how the web application would/should act in certain situations."

# --- Question 2: new list paragraph after Answer 1 ---
$r = $d.Content
$r.Find.Execute("would/should act in certain situations.")
$r.Collapse(0)
$r.InsertParagraphAfter()
$q2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$q2.Style = "List Paragraph"
$q2.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
$q2.Range.Text = "What is the difference between web server controls and HTML controls?"

# --- Answer 2: new plain paragraph after Question 2 ---
$r = $d.Content
$r.Find.Execute("What is the difference between web server controls and HTML controls?")
$r.Collapse(0)
$r.InsertParagraphAfter()
$ans2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$ans2.Range.ListFormat.RemoveNumbers()
$ans2.Style = "Normal"
$ans2.Range.Text = "Web Server Controls are a controls that derive from the System.Web.UI.WebControls base class and can be programmed to behave exactly like HTML controls, the difference is however is that they are executed on the server-side and can detect the target browser" + [char]0x2019 + "s capabilities and render themselves accordingly."

# --- Question 3: new list paragraph after Answer 2 ---
$r = $d.Content
$r.Find.Execute("render themselves accordingly.")
$r.Collapse(0)
$r.InsertParagraphAfter()
$q3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$q3.Style = "List Paragraph"
$q3.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)
$q3.Range.Text = "What is the role of property IsPostBack and how should it be used?"

# --- Answer 3: new plain paragraph after Question 3; keeps the bookmark ---
$r = $d.Content
$r.Find.Execute("What is the role of property IsPostBack and how should it be used?")
$r.Collapse(0)
$r.InsertParagraphAfter()
$ans3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$ans3.Range.ListFormat.RemoveNumbers()
$ans3.Style = "Normal"
$ans3.Range.Text = "The role of IsPostBack is to determine whether the page is being rendered for the first time or is being loaded in a response to a postback."
